$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 26.86490966666667
$ws.Range("H2").Value = 80.594729
$ws.Range("I2").Value = 0.1447302967754861
$ws.Range("J2").Value = 0.1447302967754861
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3389413333333333
$ws.Range("N2").Value = 1.016824
$ws.Range("Q2").Value = 9.105628302299555
$ws.Range("R2").Value = 81.95065472069599
$ws.Range("S2").Value = 0.1447302967754861
$ws.Range("T2").Value = 0.1447302967754861

# Row 3
$ws.Range("I3").Value = 0.4077186109324291
$ws.Range("J3").Value = 0.4077186109324292
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3389413333333333
$ws.Range("N3").Value = 1.016824
$ws.Range("Q3").Value = 25.65139577402844
$ws.Range("R3").Value = 230.862561966256
$ws.Range("S3").Value = 0.4077186109324291
$ws.Range("T3").Value = 0.4077186109324292

# Row 4
$ws.Range("G4").Value = 14.45399366666666
$ws.Range("H4").Value = 43.36198099999999
$ws.Range("I4").Value = 0.07786852138807973
$ws.Range("J4").Value = 0.07786852138807973
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3389413333333333
$ws.Range("N4").Value = 1.016824
$ws.Range("Q4").Value = 4.899055885371554
$ws.Range("R4").Value = 44.09150296834399
$ws.Range("S4").Value = 0.07786852138807973
$ws.Range("T4").Value = 0.07786852138807973

# Row 5
$ws.Range("G5").Value = 68.62066266666666
$ws.Range("H5").Value = 205.861988
$ws.Range("I5").Value = 0.369682570904005
$ws.Range("J5").Value = 0.369682570904005
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3389413333333333
$ws.Range("N5").Value = 1.016824
$ws.Range("Q5").Value = 23.25837889845689
$ws.Range("R5").Value = 209.325410086112
$ws.Range("S5").Value = 0.369682570904005
$ws.Range("T5").Value = 0.369682570904005
